$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps values as text so ambiguous numeric-looking
# strings (e.g. "66.152.47", "0.544") are not reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '66.152.47'
$ws.Range("E2").Value = '  -1.34%  '

# Row 3
$ws.Range("D3").Value = '3.204.98'
$ws.Range("E3").Value = '  -0.39%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").Value = '607.52'
$ws.Range("E5").Value = '  +0.40%  '

# Row 6
$ws.Range("D6").Value = '154.71'
$ws.Range("E6").Value = '  -1.90%  '

# Row 7
$ws.Range("E7").Value = '  +0.10%  '

# Row 8
$ws.Range("D8").Value = '3.205.89'
$ws.Range("E8").Value = '  -0.32%  '

# Row 9
$ws.Range("D9").Value = '0.544'
$ws.Range("E9").Value = '  -2.46%  '

# Row 10
$ws.Range("D10").Value = '0.159'
$ws.Range("E10").Value = '  -1.91%  '

# Row 11
$ws.Range("D11").Value = '5.69'
$ws.Range("E11").Value = '  -4.14%  '

# Row 12
$ws.Range("D12").Value = '0.498'
$ws.Range("E12").Value = '  -4.24%  '

# Row 13
$ws.Range("D13").Value = '0.0000265'
$ws.Range("E13").Value = '  -2.04%  '

# Row 14
$ws.Range("D14").Value = '38.11'
$ws.Range("E14").Value = '  -3.74%  '

# Row 15
$ws.Range("D15").Value = '3.741.06'
$ws.Range("E15").Value = '  -0.13%  '

# Row 16
$ws.Range("D16").Value = '66.335.16'
$ws.Range("E16").Value = '  -0.94%  '

# Row 17
$ws.Range("D17").Value = '3.212.40'
$ws.Range("E17").Value = '  -0.04%  '

# Row 18
$ws.Range("D18").Value = '7.22'
$ws.Range("E18").Value = '  -4.17%  '

# Row 19
$ws.Range("E19").Value = '  +0.78%  '

# Row 20
$ws.Range("D20").Value = '503.46'
$ws.Range("E20").Value = '  -4.57%  '

# Row 21
$ws.Range("D21").Value = '15.15'
$ws.Range("E21").Value = '  -2.92%  '

# Row 22
$ws.Range("D22").Value = '0.725'
$ws.Range("E22").Value = '  -3.15%  '

# Row 23
$ws.Range("D23").Value = '7.96'
$ws.Range("E23").Value = '  -3.90%  '

# Row 24
$ws.Range("D24").Value = '14.49'
$ws.Range("E24").Value = '  -4.37%  '

# Row 25
$ws.Range("D25").Value = '84.68'
$ws.Range("E25").Value = '  -1.58%  '

# Row 26
$ws.Range("B26").Value = 'Hedera'
$ws.Range("C26").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D26").Value = '0.157'
$ws.Range("E26").Value = '  +72.18%  '

# Row 27
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.20%  '

# Row 28
$ws.Range("D28").Value = '2.99'
$ws.Range("E28").Value = '  -1.42%  '

# Row 29
$ws.Range("D29").Value = '8.91'
$ws.Range("E29").Value = '  -4.67%  '

# Row 30
$ws.Range("D30").Value = '2.34'
$ws.Range("E30").Value = '  -3.16%  '

# Row 31
$ws.Range("D31").Value = '6.90'
$ws.Range("E31").Value = '  -2.36%  '

# Row 32
$ws.Range("D32").Value = '2.87'
$ws.Range("E32").Value = '  -3.75%  '

# Row 33
$ws.Range("D33").Value = '28.01'
$ws.Range("E33").Value = '  -1.70%  '

# Row 34
$ws.Range("E34").Value = '  +0.06%  '

# Row 35
$ws.Range("E35").Value = '  -5.54%  '

# Row 36
$ws.Range("D36").Value = '6.35'
$ws.Range("E36").Value = '  -4.12%  '

# Row 37
$ws.Range("D37").Value = '55.31'
$ws.Range("E37").Value = '  +0.14%  '

# Row 38
$ws.Range("D38").Value = '495.69'
$ws.Range("E38").Value = '  -5.49%  '

# Row 39
$ws.Range("D39").Value = '0.0₃0761'
$ws.Range("E39").Value = '  +10.67%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '0.129'
$ws.Range("E40").Value = '  +1.47%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0416'
$ws.Range("E41").Value = '  -3.05%  '

# Row 42
$ws.Range("D42").Value = '2.98'
$ws.Range("E42").Value = '  +1.25%  '

# Row 43
$ws.Range("D43").Value = '8.68'
$ws.Range("E43").Value = '  -3.01%  '

# Row 44
$ws.Range("D44").Value = '0.293'
$ws.Range("E44").Value = '  -3.47%  '

# Row 45
$ws.Range("D45").Value = '2.913.95'
$ws.Range("E45").Value = '  +0.06%  '

# Row 46
$ws.Range("D46").Value = '2.42'
$ws.Range("E46").Value = '  -2.65%  '

# Row 47
$ws.Range("D47").Value = '27.85'
$ws.Range("E47").Value = '  -3.53%  '

# Row 48
$ws.Range("D48").Value = '2.37'
$ws.Range("E48").Value = '  +0.19%  '

# Row 50
$ws.Range("E50").Value = '  -1.69%  '

# Row 51
$ws.Range("D51").Value = '121.61'
$ws.Range("E51").Value = '  -0.53%  '
